$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (column D) and Volume(1h) (column E) values for the
# cryptos list, as refreshed by the scheduled GitHub Actions job.
$updates = @(
    @{ Row = 2; D = "30.216.92"; E = "  -0.59%  " },
    @{ Row = 3; D = "1.862.14"; E = "  -0.66%  " },
    @{ Row = 4; D = "1.001"; E = "  +0.16%  " },
    @{ Row = 5; D = "236.24"; E = "  +0.19%  " },
    @{ Row = 6; D = "1.001"; E = "  +0.06%  " },
    @{ Row = 7; D = $null; E = "  +0.27%  " },
    @{ Row = 8; D = "0.2875"; E = "  +1.13%  " },
    @{ Row = 9; D = "0.06535"; E = "  -0.18%  " },
    @{ Row = 10; D = "21.68"; E = "  +1.32%  " },
    @{ Row = 11; D = "0.07940"; E = "  +0.31%  " },
    @{ Row = 12; D = "97.77"; E = "  +0.45%  " },
    @{ Row = 13; D = "1.869.25"; E = "  -0.32%  " },
    @{ Row = 14; D = "5.168"; E = "  +0.13%  " },
    @{ Row = 15; D = "0.6800"; E = "  +0.72%  " },
    @{ Row = 16; D = "267.27"; E = "  -5.67%  " },
    @{ Row = 17; D = "30.221.41"; E = "  -0.59%  " },
    @{ Row = 18; D = $null; E = "  +7.89%  " },
    @{ Row = 19; D = "1.001"; E = "  +0.11%  " },
    @{ Row = 20; D = "0.000007398"; E = "  +1.45%  " },
    @{ Row = 21; D = "2.114.19"; E = "  -0.62%  " },
    @{ Row = 22; D = "5.310"; E = "  -4.50%  " },
    @{ Row = 23; D = $null; E = "  +0.17%  " },
    @{ Row = 24; D = "6.168"; E = "  -0.84%  " },
    @{ Row = 25; D = "167.50"; E = "  +1.41%  " },
    @{ Row = 26; D = "9.202"; E = "  -1.23%  " },
    @{ Row = 27; D = "18.84"; E = "  -1.51%  " },
    @{ Row = 28; D = "1.954"; E = "  +0.49%  " },
    @{ Row = 29; D = $null; E = "  +1.85%  " },
    @{ Row = 30; D = "0.09813"; E = "  +1.20%  " },
    @{ Row = 31; D = "4.374"; E = "  -1.55%  " },
    @{ Row = 32; D = "1.469"; E = "  -0.47%  " },
    @{ Row = 33; D = "4.058"; E = "  -1.46%  " },
    @{ Row = 34; D = "0.04690"; E = "  -0.50%  " },
    @{ Row = 35; D = "1.130"; E = "  +0.39%  " },
    @{ Row = 36; D = "0.7008"; E = "  -0.68%  " },
    @{ Row = 37; D = "2.708"; E = "  -0.38%  " },
    @{ Row = 38; D = "0.01872"; E = "  +0.37%  " },
    @{ Row = 39; D = "2.620"; E = "  +3.24%  " },
    @{ Row = 40; D = "6.248"; E = "  -2.03%  " },
    @{ Row = 41; D = "74.43"; E = "  +0.20%  " },
    @{ Row = 42; D = "1.935"; E = "  -1.33%  " },
    @{ Row = 44; D = "0.4161"; E = "  -0.97%  " },
    @{ Row = 45; D = "0.9996"; E = "  -0.02%  " },
    @{ Row = 46; D = "103.08"; E = "  -0.76%  " },
    @{ Row = 47; D = "952.84"; E = "  +1.86%  " },
    @{ Row = 48; D = "7.152"; E = "  -1.43%  " },
    @{ Row = 49; D = "9.206"; E = "  -1.26%  " },
    @{ Row = 50; D = "34.12"; E = "  -0.28%  " },
    @{ Row = 51; D = "0.05667"; E = "  +0.66%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        # Force the assignment to be kept as text (these prices use "."
        # as thousands separators, e.g. "1.862.14", and some values such
        # as "1.001" would otherwise be auto-coerced into a number),
        # then drop back to the default "Normal" style so no extra
        # formatting is left behind on the cell.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($null -ne $u.E) {
        $cell = $ws.Range("E" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
        $cell.Style = "Normal"
    }
}
